$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old index column (A), shifting col1..col4 (B:E) left into A:D.
$ws.Range("A1").EntireColumn.Delete()

$ws.Range("D4").Select()
